# remove use of slurmtools on slide 22
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)
$shp = $s.Shapes.Item(6)

# Delete the first paragraph ("$ module load slurmtools") entirely,
# collapsing the text box down to just the "$ seff <job number>" line.
$tr = $shp.TextFrame.TextRange
$firstPara = $tr.Paragraphs(1, 1)
$firstPara.Delete()

# Shrink the shape to its new (shorter) height now that a line is gone.
$shp.Height = 40.054
